$wb = $excel.ActiveWorkbook

# --- New report identifiers (this handoff run) -----------------------------
$oldGuid = "3625d91a-21e7-49c5-b1b5-33c6a6d0d442"
$newGuid = "c6159970-c0da-4760-9c8e-2a9162e7e16e"

$oldSha = "f9caf711a3bf7ce24dc34b24edfc6206c9756dc0"
$newSha = "28d74f79a3fee63c8bffb2271984207348cacaa9"

# Hyperlink targets keep pointing at the same (already-published) commit/file,
# only the visible/display text is refreshed to the new guid.
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e10c3e0940e9e829a17fe2965af7b0592d7961cf/e2e/$oldGuid.md"

# Helper: refresh the hyperlink on a cell (delete + re-add is the only way to
# update a hyperlink's display text without leaving a stale duplicate behind)
# and then restore the underline/blue-link look the cell had before.
function Update-HandoffLink($ws, $cellAddr, $displayText) {
    $cell = $ws.Range($cellAddr)
    $cell.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($cell, $hyperlinkUrl, "", "", $displayText) | Out-Null
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Underline = $true
    $cell.Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
Update-HandoffLink $wsOverview "B2" "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-23 12:58:27"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-HandoffLink $wsZh "A2" "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newSha.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-23 12:58:22"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
Update-HandoffLink $wsDe "A2" "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newSha.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-23 12:58:27"
